$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update URL column (B) for rows 2-6: south-africa -> puerto-rico
$ws.Range("B2:B6").Value = "https://www.varoom.com/all/puerto-rico"

# Row 2: price -> rating, values -> New
$ws.Range("E2").Value = "rating"
$ws.Range("F2").Value = "New"
$ws.Range("G2").Value = "New"
$ws.Range("H2").Value = "New"

# Row 3: rating -> number_of_reviews, values -> New
$ws.Range("E3").Value = "number_of_reviews"
$ws.Range("F3").Value = "New"
$ws.Range("G3").Value = "New"
$ws.Range("H3").Value = "New"

# Row 4: number_of_reviews -> property_type, values -> Apartment
$ws.Range("E4").Value = "property_type"
$ws.Range("F4").Value = "Apartment"
$ws.Range("G4").Value = "Apartment"
$ws.Range("H4").Value = "Apartment"

# Row 5: property_type -> title, values -> new title text
$ws.Range("E5").Value = "title"
$ws.Range("F5").Value = "Newly Furnished 3 Bedroom Apartment 1 Bath in Hato Rey San Juan Puerto Rico"
$ws.Range("G5").Value = "Newly Furnished 3 Bedroom Apartment 1 Bath in Hato Rey San Juan Puerto Rico"
$ws.Range("H5").Value = "Newly Furnished 3 Bedroom Apartment 1 Bath in Hato Rey San Juan Puerto Rico"

# Row 6: title -> price, values -> new price, Passed flag flips to TRUE
$ws.Range("E6").Value = "price"
$ws.Range("F6").Value = "৳14,590"
$ws.Range("G6").Value = "৳14,590"
$ws.Range("H6").Value = "৳14,590"
$ws.Range("I6").Value = $true
